# Requirement Specification workbook update
# - Rename the "UI" sub-acronym/category to "Common" ("U"/"UI"/"UI 관련" -> "C"/"Common"/"광역 관련")
# - Replace the old "FR-U-A-xx" (UI/Animation) entries with the new Mouse Event Manager
#   "FR-C-xx" entries
# - Add two new "FR-C-A-xx" rows describing the animation-related mouse work
# - Rename "FR-U-T-xx" codes to "FR-C-T-xx"
# - Mark a few Implementation status cells as complete (green)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constant (Excel VBA/COM Color is a BGR integer, not RGB)
$colorGreen  = 5287936   # theme fill used for "완료" (Complete) -> RGB 00B050

# --- Acronym legend row (row 6): "U" / "UI" / "UI 관련" -> "C" / "Common" / "광역 관련" ---
$ws.Range("E6").Value = "C"
$ws.Range("F6").Value = "Common"
$ws.Range("G6").Value = "광역 관련"

# --- Row 23: old "UI / A" (FR-U-A-01, Fade in/out) -> new "Common" mouse-event entry ---
$ws.Range("D23").Value = "C"
$ws.Range("F23").ClearContents()
$ws.Range("H23").Value = "FR-C-01"
$ws.Range("J23").Value = "마우스 이벤트 메니저"

# --- Row 24: old FR-U-A-02 -> FR-C-02 ---
$ws.Range("H24").Value = "FR-C-02"
$ws.Range("J24").Value = "마우스 이벤트 메니저에 이벤트 전달"

# --- Row 28: new "A" (Animation) sub-section, first entry (previously blank row) ---
$ws.Range("F28").Value = "A"
$ws.Range("H28").Value = "FR-C-A-01"
$ws.Range("J28").Value = "게임 오브젝트 Fade in/out"
$ws.Range("S28").Interior.Color = $colorGreen

# --- Row 29: new "A" (Animation) sub-section, second entry (previously blank row) ---
$ws.Range("H29").Value = "FR-C-A-02"
$ws.Range("J29").Value = "게임 오브젝트 수직/수평 으로 입력한 거리만큼 일정/가속 속도로 이동"
$ws.Range("S29").Interior.Color = $colorGreen

# --- Row 37/38: rename Notification codes FR-U-T-0x -> FR-C-T-0x ---
$ws.Range("H37").Value = "FR-C-T-01"
$ws.Range("H38").Value = "FR-C-T-02"

# --- Implementation status: mark FR-S-01 (row 13) as complete (was waiting) ---
$ws.Range("S13").Interior.Color = $colorGreen

# --- Leave the selection where the author's session ended up ---
[void]$ws.Range("W20").Select()
